# importar_cotizaciones.xlsx update:
#   - add a new "MONTO MMOO" column (G) next to the existing quote columns
#   - clear out the stray formatted-but-empty cells in E:G (rows 2-13) and
#     the three fully-empty extra rows (11-13) that had no real data
#   - give the new helper cell G6 an underlined font (matches the author's
#     "close session" note - an underlined link-like cell)
#   - move the active selection to G6 and size the new column to fit its header

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) New header in G1
$ws.Range("G1").Value = "MONTO MMOO"

# Keep the whole header row's alignment consistent (centered), same as the
# existing A1:F1 headers already were.
$ws.Range("A1:G1").HorizontalAlignment = -4108   # xlCenter

# 2) Remove the left-over formatted-only cells in columns E:G for the data
#    rows, and drop rows 11-13 completely (they only ever held empty,
#    formatted cells, no real data).
$ws.Range("E2:G13").Clear()
$ws.Range("A11:G13").Clear()

# 3) New styled (underlined) helper cell at G6
$ws.Range("G6").Font.Underline = $true

# 4) Resize column G to fit the new header text
$ws.Columns("G").ColumnWidth = 13.833333333333334

# 5) Update the active selection to G6
$ws.Range("G6").Select() | Out-Null
